# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Titan_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 366.66666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 366.66666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 366.66666
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -706.66666

$ws.Range("H17").Value = 622.61176
$ws.Range("J17").Value = 627.9759
$ws.Range("L17").Value = 1883.9277
$ws.Range("N17").Value = -2219.9277

$ws.Range("H31").Value = 5692
$ws.Range("I31").Value = 5692
$ws.Range("K31").Value = 17076
$ws.Range("M31").Value = -16846

$ws.Range("H106").Value = 7009483.5
$ws.Range("I106").Value = 8010459.5
$ws.Range("J106").Value = 2650
$ws.Range("K106").Value = 8010459.5
$ws.Range("L106").Value = 2650
$ws.Range("M106").Value = -8009828.5
$ws.Range("N106").Value = -3912

$ws.Range("H137").Value = 41668332
$ws.Range("I137").Value = 58824610
$ws.Range("J137").Value = 3081.5715
$ws.Range("K137").Value = 176473830
$ws.Range("L137").Value = 9244.7145
$ws.Range("M137").Value = -176471280
$ws.Range("N137").Value = -14344.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2123.6875
$ws.Range("I2").Value = 1977.0714
$ws.Range("J2").Value = 3150
$ws.Range("K2").Value = 1977.0714
$ws.Range("L2").Value = 3150
$ws.Range("M2").Value = -1864.0714
$ws.Range("N2").Value = -3376

$ws.Range("H45").Value = 1516.6666
$ws.Range("I45").Value = 1350
$ws.Range("J45").Value = 1600
$ws.Range("K45").Value = 1350
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = -973
$ws.Range("N45").Value = -2354

$ws.Range("H61").Value = 3418.2
$ws.Range("I61").Value = 2053.4285
$ws.Range("K61").Value = 2053.4285
$ws.Range("M61").Value = -1841.4285

$ws.Range("H74").Value = 4965.4116
$ws.Range("I74").Value = 1452.44
$ws.Range("J74").Value = 14723.667
$ws.Range("K74").Value = 1452.44
$ws.Range("L74").Value = 14723.667
$ws.Range("M74").Value = -578.4400000000001
$ws.Range("N74").Value = -16471.667

$ws.Range("H77").Value = 4965.4116
$ws.Range("I77").Value = 1452.44
$ws.Range("J77").Value = 14723.667
$ws.Range("K77").Value = 7262.200000000001
$ws.Range("L77").Value = 73618.33499999999
$ws.Range("M77").Value = -2894.200000000001
$ws.Range("N77").Value = -82354.33499999999

$ws.Range("H97").Value = 17549676
$ws.Range("I97").Value = 22229192
$ws.Range("J97").Value = 1486.75
$ws.Range("K97").Value = 22229192
$ws.Range("L97").Value = 1486.75
$ws.Range("M97").Value = -22228696
$ws.Range("N97").Value = -2478.75

$ws.Range("H116").Value = 2123.6875
$ws.Range("I116").Value = 1977.0714
$ws.Range("J116").Value = 3150
$ws.Range("K116").Value = 1977.0714
$ws.Range("L116").Value = 3150
$ws.Range("M116").Value = 316.9286
$ws.Range("N116").Value = -7738

$ws.Range("H122").Value = 2026.7273
$ws.Range("I122").Value = 1703
$ws.Range("J122").Value = 2211.7144
$ws.Range("K122").Value = 5109
$ws.Range("L122").Value = 6635.1432
$ws.Range("M122").Value = -2659
$ws.Range("N122").Value = -11535.1432

$ws.Range("H132").Value = 3562.0715
$ws.Range("I132").Value = 2895.647
$ws.Range("J132").Value = 4592
$ws.Range("K132").Value = 8686.940999999999
$ws.Range("L132").Value = 13776
$ws.Range("M132").Value = -6156.940999999999
$ws.Range("N132").Value = -18836

$ws.Range("H136").Value = 3418.2
$ws.Range("I136").Value = 2053.4285
$ws.Range("K136").Value = 6160.2855
$ws.Range("M136").Value = -3610.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2123.6875
$ws.Range("I3").Value = 1977.0714
$ws.Range("J3").Value = 3150
$ws.Range("K3").Value = 1977.0714
$ws.Range("L3").Value = 3150
$ws.Range("M3").Value = -1863.0714
$ws.Range("N3").Value = -3378

$ws.Range("H105").Value = 3056.6
$ws.Range("I105").Value = 2953.4211
$ws.Range("J105").Value = 3383.3333
$ws.Range("K105").Value = 2953.4211
$ws.Range("L105").Value = 3383.3333
$ws.Range("M105").Value = -1206.4211
$ws.Range("N105").Value = -6877.3333

$ws.Range("H134").Value = 3531.111
$ws.Range("I134").Value = 2456.05
$ws.Range("K134").Value = 7368.150000000001
$ws.Range("M134").Value = -4833.150000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1521.75
$ws.Range("I31").Value = 1167.7142
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1167.7142
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -872.7141999999999
$ws.Range("N31").Value = -4590

$ws.Range("H32").Value = 16803.2
$ws.Range("I32").Value = 16803.2
$ws.Range("K32").Value = 16803.2
$ws.Range("M32").Value = -16487.2

$ws.Range("H34").Value = 1521.75
$ws.Range("I34").Value = 1167.7142
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1167.7142
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -965.7141999999999
$ws.Range("N34").Value = -4404

$ws.Range("H58").Value = 2237.9412
$ws.Range("I58").Value = 1235.762
$ws.Range("K58").Value = 1235.762
$ws.Range("M58").Value = -1032.762

$ws.Range("H105").Value = 471.4375
$ws.Range("I105").Value = 476.33334
$ws.Range("J105").Value = 465.14285
$ws.Range("K105").Value = 476.33334
$ws.Range("L105").Value = 465.14285
$ws.Range("M105").Value = 1270.66666
$ws.Range("N105").Value = -3959.14285

$ws.Range("H132").Value = 3384.52
$ws.Range("I132").Value = 2519
$ws.Range("J132").Value = 4682.8
$ws.Range("K132").Value = 7557
$ws.Range("L132").Value = 14048.4
$ws.Range("M132").Value = -5027
$ws.Range("N132").Value = -19108.4

$ws.Range("H134").Value = 3416.111
$ws.Range("I134").Value = 1720.25
$ws.Range("J134").Value = 4772.8
$ws.Range("K134").Value = 5160.75
$ws.Range("L134").Value = 14318.4
$ws.Range("M134").Value = -2625.75
$ws.Range("N134").Value = -19388.4

$ws.Range("H136").Value = 2237.9412
$ws.Range("I136").Value = 1235.762
$ws.Range("K136").Value = 3707.286
$ws.Range("M136").Value = -1157.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 889.4706
$ws.Range("I5").Value = 651.4286
$ws.Range("K5").Value = 1954.2858
$ws.Range("M5").Value = -1842.2858

$ws.Range("H135").Value = 889.4706
$ws.Range("I135").Value = 651.4286
$ws.Range("K135").Value = 5862.8574
$ws.Range("M135").Value = -3327.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1859.5
$ws.Range("I102").Value = 1512.875
$ws.Range("J102").Value = 2321.6667
$ws.Range("K102").Value = 1512.875
$ws.Range("L102").Value = 2321.6667
$ws.Range("M102").Value = 109.125
$ws.Range("N102").Value = -5565.6667

$ws.Range("H132").Value = 2883.451
$ws.Range("I132").Value = 2339.6667
$ws.Range("J132").Value = 3660.2856
$ws.Range("K132").Value = 7019.000100000001
$ws.Range("L132").Value = 10980.8568
$ws.Range("M132").Value = -4489.000100000001
$ws.Range("N132").Value = -16040.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 3693.7817
$ws.Range("I132").Value = 2970.8215
$ws.Range("J132").Value = 4443.5186
$ws.Range("K132").Value = 8912.4645
$ws.Range("L132").Value = 13330.5558
$ws.Range("M132").Value = -6382.4645
$ws.Range("N132").Value = -18390.5558

$ws.Range("H136").Value = 3522.6829
$ws.Range("I136").Value = 1850.931
$ws.Range("J136").Value = 7562.75
$ws.Range("K136").Value = 5552.793
$ws.Range("L136").Value = 22688.25
$ws.Range("M136").Value = -3002.793
$ws.Range("N136").Value = -27788.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1888
$ws.Range("N5").ClearContents()

$ws.Range("H20").Value = 3000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -3480

$ws.Range("H23").Value = 6790.143
$ws.Range("I23").Value = 3005
$ws.Range("J23").Value = 11837
$ws.Range("K23").Value = 3005
$ws.Range("L23").Value = 11837
$ws.Range("M23").Value = -2776
$ws.Range("N23").Value = -12295

$ws.Range("H74").Value = 13531
$ws.Range("J74").Value = 13723.4
$ws.Range("L74").Value = 13723.4
$ws.Range("N74").Value = -15595.4

$ws.Range("H77").Value = 13531
$ws.Range("J77").Value = 13723.4
$ws.Range("L77").Value = 41170.2
$ws.Range("N77").Value = -50530.2

$ws.Range("H107").Value = 2137140
$ws.Range("I107").Value = 3086830.5
$ws.Range("J107").Value = 336.25
$ws.Range("K107").Value = 9260491.5
$ws.Range("L107").Value = 1008.75
$ws.Range("M107").Value = -9258571.5
$ws.Range("N107").Value = -4848.75

$ws.Range("H108").Value = 49990
$ws.Range("J108").Value = 49990
$ws.Range("L108").Value = 49990
$ws.Range("N108").Value = -57670

$ws.Range("H122").Value = 40250.27
$ws.Range("J122").Value = 1908.4
$ws.Range("L122").Value = 5725.200000000001
$ws.Range("N122").Value = -10625.2

$ws.Range("H123").Value = 24869.084
$ws.Range("J123").Value = 24869.084
$ws.Range("L123").Value = 24869.084
$ws.Range("N123").Value = -34669.084

$ws.Range("H132").Value = 23813074
$ws.Range("I132").Value = 38464590
$ws.Range("K132").Value = 115393770
$ws.Range("M132").Value = -115391240

